$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expense")

# Use a scratch cell far away to produce literal text values (via formula + paste-values)
# so Excel's automatic date/number recognition doesn't convert "01/01/2109" into a date
# serial and doesn't leave a residual number-format style on the target cells.
$scratch = $ws.Cells.Item(100, 1)

function Set-LiteralText($cell, $text) {
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
}

Set-LiteralText $ws.Cells.Item(5,1) "01/01/2109"
Set-LiteralText $ws.Cells.Item(5,2) "Food"
$ws.Cells.Item(5,3).NumberFormat = "$#,##0.00"
$ws.Cells.Item(5,3).Value = 9000

Set-LiteralText $ws.Cells.Item(6,1) "01/01/2109"
Set-LiteralText $ws.Cells.Item(6,2) "Food"
$ws.Cells.Item(6,3).Value = 9000

Set-LiteralText $ws.Cells.Item(7,1) "01/01/2109"
Set-LiteralText $ws.Cells.Item(7,2) "Food"
$ws.Cells.Item(7,3).Value = 9000

$scratch.ClearContents()
